$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (DOI, Mistake,
# Fig ID, Journal, ...) from A:J to B:K.
$ws.Columns.Item(1).Insert()

# Fig Index = "<DOI>_<Mistake>_fig<N>", where N increments per repeated DOI
# within the same paper -- pairs each annotation row with its figure so it
# can later be joined against the distortion-metric table.
$figIndexValues = @{
    13 = "10.1002:btm2.10634_log_fig1"
    12 = "10.1002:btm2.10602_log_fig1"
    11 = "10.1002:btm2.10594_log_fig4"
    10 = "10.1002:btm2.10594_log_fig3"
    9  = "10.1002:btm2.10594_log_fig2"
    8  = "10.1002:btm2.10594_log_fig1"
    7  = "10.1002:btm2.10538_log_fig2"
    6  = "10.1002:btm2.10538_log_fig1"
    5  = "10.1002:btm2.10514_log_fig2"
    4  = "10.1002:btm2.10514_log_fig1"
    3  = "10.1002:btm2.10487_log_fig2"
    2  = "10.1002:btm2.10487_log_fig1"
}

foreach ($row in 13..2) {
    $ws.Cells.Item($row, 1).Value = $figIndexValues[$row]
}

# Header for the new "Fig Index" column (typed in last, matching the style of
# the adjacent "DOI" header cell).
$ws.Cells.Item(1, 1).Value = "Fig Index"

$ws.Range("B1").Copy()
$ws.Range("A1:A13").PasteSpecial(-4122)  # xlPasteFormats

# Restore the active selection that the author left on save.
$ws.Range("B5").Select()
